# Further Data Exploration and Visualisation
#
# Week 7's task description changes focus: instead of leading with the
# model-refinement step, the plan now leads with the write-up, keeping the
# refinement as a conditional follow-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "Focus on write-up. If needed, refine the machine learning models based on the insights gained from experiments and observations. Fine-tune the models' hyperparameters to optimize their performance. Validate the models' generalizability by testing them on unseen data."

# Sheet view changes: zoom to 90%, scroll the view so column B is visible at
# the left edge, and leave the selection on B12.
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B12").Select()
